$d = $word.ActiveDocument

# Update the date heading at the top of the document
# (2024-08-26 Monday -> 2024-08-27 Tuesday)
$dateParagraph = $d.Paragraphs.Item(1)
$dateParagraph.Range.Text = "2024-08-27 Tuesday"

# New arithmetic expressions for each table cell, in row-major order
# (row 1 col 1..5, row 2 col 1..5, ... row 20 col 1..5)
$newValues = @(
    "39-1=",
    "69+15=",
    "86-1=",
    "7+82=",
    "80-26=",
    "94-0=",
    "3+67=",
    "6+44=",
    "86-25=",
    "89-17=",
    "99-80=",
    "15+77=",
    "64-29=",
    "45-3=",
    "61-34=",
    "53-17=",
    "96-41=",
    "20+4=",
    "73+2=",
    "25+66=",
    "44+54=",
    "62-12=",
    "54-27=",
    "39+24=",
    "59-0=",
    "90+1=",
    "97-17=",
    "66-60=",
    "28-2=",
    "18+20=",
    "29+60=",
    "81-51=",
    "60+26=",
    "37+30=",
    "28-10=",
    "97-11=",
    "56+12=",
    "86-4=",
    "97-0=",
    "96-34=",
    "2+36=",
    "90-4=",
    "57+23=",
    "61+38=",
    "85-11=",
    "1+68=",
    "0+84=",
    "87-82=",
    "42-13=",
    "42+3=",
    "64+13=",
    "4+9=",
    "10+71=",
    "11+65=",
    "74-33=",
    "84-5=",
    "31+50=",
    "42+25=",
    "30+50=",
    "89-60=",
    "37+20=",
    "2+39=",
    "71-18=",
    "17+44=",
    "19+56=",
    "73+0=",
    "51+6=",
    "89-1=",
    "62-22=",
    "66-10=",
    "65-53=",
    "21+71=",
    "40+3=",
    "40+18=",
    "29-13=",
    "30-7=",
    "6-1=",
    "97-59=",
    "8+56=",
    "61-25=",
    "9+63=",
    "72+20=",
    "74-57=",
    "14+53=",
    "76+13=",
    "20-18=",
    "57-22=",
    "79-70=",
    "97-79=",
    "62-53=",
    "79-25=",
    "49-8=",
    "9+20=",
    "68-48=",
    "49-42=",
    "42+35=",
    "13+69=",
    "53-12=",
    "2+33=",
    "73-52="
)

$tbl = $d.Tables.Item(1)
$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        # Assign directly to the cell range text (rather than Find/Replace)
        # to keep each edit precisely scoped to its own cell.
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated " + $idx + " cells")
